$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.256.47"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.645.62"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.85%  "
$ws.Range("D9").Value = "2.645.04"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  +1.57%  "
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "3.129.77"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").Value = "72.217.46"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").Value = "2.707.36"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("D28").Value = "2.783.32"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "0.0₃0956"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "500.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -7.02%  "
$ws.Range("E43").Value = "  -3.48%  "
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.67%  "
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  -1.28%  "
